$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$ws1.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"
$ws1.Range("B12").Value = "Expandable categorisation of consumer/requestor of a consent enforcement request related to a specific consent component (e.g. policy or module) "
